$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Relocate the "_GoBack" bookmark: remove it from its old spot (right
#    after "...ensure security of cloud applications.") - it will be
#    re-added at the end of the newly inserted paragraph below.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2. Replace the lone "space only" paragraph that sits between the blank
#    paragraph following "Critique - Sogra" and the "References" heading
#    with the new paragraph of text (purpose/overview of the paper),
#    re-creating the mid-sentence lastRenderedPageBreak marker and putting
#    the "_GoBack" bookmark back at the end of it.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$critiqueIdx = $full.IndexOf("Critique")
$prefix = $full.Substring(0, $critiqueIdx + 19)
$parMarks = ($prefix.ToCharArray() | Where-Object { $_ -eq [char]13 }).Count
$target = $d.Paragraphs.Item($parMarks + 1)
$targetRange = $target.Range

$part1 = "The purpose of this research paper is to present a solution to the challenges faced in the development of secure and dependable cloud applications. This paper has a well-defined purpose and problem statement which have been communicated clearly to the readers. The paper follows a structured format of first defining the problem and then defining the concepts required to understand the solution being proposed. The paper has been structured to provide the readers with an overview of what the paper contains: purpose "
$part2 = "of the paper, problem statement and the proposed solution. This is followed by a description of the phases involved in the Software Development Life Cycle (SDLC) and the Data Security Life Cycle (DSLC). The paper then introduces the SaaS Security Life Cycle which is a combination of the two development life cycles SDLC and DSLC). The relationship between the three life cycles have been clearly demonstrated through a well-structured diagram. The paper then highlights the necessity of following the SaaS Security Life Cycle (SSLC) for the development of dependable and secure cloud application. The challenges faced have been split across the different phases of SSLC. Each phase of SSLC is linked back to its equivalent phases in the SDLC and DSLC life cycles defining the "

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">' + $part1 + '</w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">' + $part2 + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$targetRange.InsertXML($xmlFrag)
